$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows for the "Resolving-Mac" sending cluster (rows 6-9);
# the surviving data is just the "FAPs" sending-cluster rows (2-5).
$ws.Range("A6:T9").EntireRow.Delete()

# Update the recomputed TPM-derived statistics for the remaining rows.
$ws.Range("I2").Value() = 1
$ws.Range("J2").Value() = 1
$ws.Range("M2").Value() = 1.376848666666667
$ws.Range("N2").Value() = 4.130546
$ws.Range("O2").Value() = 0.1003061556015877
$ws.Range("P2").Value() = 0.1003061556015877
$ws.Range("Q2").Value() = 0.1612119977331111
$ws.Range("R2").Value() = 1.450907979598
$ws.Range("S2").Value() = 0.1003061556015877
$ws.Range("T2").Value() = 0.1003061556015877

$ws.Range("I3").Value() = 1
$ws.Range("J3").Value() = 1
$ws.Range("M3").Value() = 2.848096333333333
$ws.Range("N3").Value() = 8.544288999999999
$ws.Range("O3").Value() = 0.207489465542554
$ws.Range("P3").Value() = 0.207489465542554
$ws.Range("Q3").Value() = 0.3334769541118888
$ws.Range("S3").Value() = 0.207489465542554
$ws.Range("T3").Value() = 0.207489465542554

$ws.Range("I4").Value() = 1
$ws.Range("J4").Value() = 1
$ws.Range("M4").Value() = 7.562766000000001
$ws.Range("N4").Value() = 22.688298
$ws.Range("O4").Value() = 0.5509624997574636
$ws.Range("P4").Value() = 0.5509624997574636
$ws.Range("Q4").Value() = 0.885506624486
$ws.Range("R4").Value() = 7.969559620374001
$ws.Range("S4").Value() = 0.5509624997574636
$ws.Range("T4").Value() = 0.5509624997574636

$ws.Range("I5").Value() = 1
$ws.Range("J5").Value() = 1
$ws.Range("M5").Value() = 1.938751333333333
$ws.Range("N5").Value() = 5.816254000000001
$ws.Range("O5").Value() = 0.1412418790983945
$ws.Range("P5").Value() = 0.1412418790983945
$ws.Range("Q5").Value() = 0.2270038698668889
$ws.Range("R5").Value() = 2.043034828802
$ws.Range("S5").Value() = 0.1412418790983945
$ws.Range("T5").Value() = 0.1412418790983945
